# Update values in column C for specific rows as described by the diff
# (data imputation algorithm result update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -12.223
    6   = -12.318
    7   = -12.909
    8   = -12.606
    16  = -13.164
    20  = -12.205
    21  = -12.725
    28  = -12.848
    29  = -11.736
    30  = -12.351
    32  = -12.95
    40  = -12.446
    46  = -13.827
    51  = -11.187
    52  = -11.614
    57  = -13.632
    59  = -13.18
    62  = -13.555
    66  = -10.836
    73  = -12.908
    74  = -11.982
    77  = -12.466
    92  = -11.104
    100 = -12.747
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
